$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("faq")

$oldUrl = "https://www.oeaw.ac.at/fileadmin/subsites/Institute/VID/PDF/Publications/Working_Papers/WP2019_02.pdf"
$newUrl = "https://www.oeaw.ac.at/fileadmin/subsites/Institute/VID/IMG/Publications/Working_Papers/WP2019_02.pdf"

# C17: "Which education categories are available?" answer references the Speringer et al., 2019 working paper.
$c17 = $ws.Range("C17")
$c17.Value = $c17.Value().Replace($oldUrl, $newUrl)

# C26: "Why is the reconstructed data from 1970 to 2010 different from Version 1.2?" answer
# also references the same working paper.
$c26 = $ws.Range("C26")
$c26.Value = $c26.Value().Replace($oldUrl, $newUrl)
